# "Creating setup file.xlsx" - add Signal name/left/top config, a piloting
# block on the Variables sheet, and rewire the formulas that build the
# config string so the signal data is folded in as well.

$wb = $excel.ActiveWorkbook
$wsParam = $wb.Worksheets.Item("Paramétrage")
$wsVar   = $wb.Worksheets.Item("Variables")

# ---------------------------------------------------------------------
# 1) Paramétrage: append the "Signal" table (rows 25-33), mirroring the
#    existing "Aiguillages" table (rows 16-24).
# ---------------------------------------------------------------------

# Header row (row 25) - same look as row 16.
$wsParam.Rows("25:25").WrapText = $true

$wsParam.Range("B25").Value = "Name (use only ansi caracters)"
$wsParam.Range("B25").Font.Bold = $true
$wsParam.Range("B25").HorizontalAlignment = -4131   # xlLeft

$wsParam.Range("C25").Value = "left"
$wsParam.Range("C25").Font.Bold = $true
$wsParam.Range("C25").HorizontalAlignment = -4108   # xlCenter

$wsParam.Range("D25").Value = "top"
$wsParam.Range("D25").Font.Bold = $true
$wsParam.Range("D25").HorizontalAlignment = -4108   # xlCenter

$wsParam.Rows("25:25").RowHeight = 45

# Data rows 26-33 (Signal 1..8 / names / left / top)
$signalLabels = @("Signal 1","Signal 2","Signal 3","Signal 4","Signal 5","Signal 6","Signal 7","Signal 8")
for ($i = 0; $i -lt 8; $i++) {
    $r = 26 + $i
    $wsParam.Range("A$r").Value = $signalLabels[$i]
    $wsParam.Range("A$r").Font.Bold = $true
}

$wsParam.Range("B26").Value = "cote gauceh"
$wsParam.Range("B27").Value = "milieu"
$wsParam.Range("B28").Value = "haut"
$wsParam.Range("B29").Value = "bas"
$wsParam.Range("B30").Value = "par ici"
$wsParam.Range("B31").Value = "par la"
$wsParam.Range("B32").Value = "par la bas"
$wsParam.Range("B33").Value = "tut en haut"

$leftVals = @(200,120,279,250,130,33,127,536)
$topVals  = @(550,10,523,160,160,645,457,546)
for ($i = 0; $i -lt 8; $i++) {
    $r = 26 + $i
    $wsParam.Range("C$r").Value = $leftVals[$i]
    $wsParam.Range("C$r").HorizontalAlignment = -4108   # xlCenter
    $wsParam.Range("D$r").Value = $topVals[$i]
    $wsParam.Range("D$r").HorizontalAlignment = -4108   # xlCenter
}

# ---------------------------------------------------------------------
# 2) Variables: insert a new row 11 ("Name Signal" / "sna") - this is
#    what shifts every row below it (old 12..29) down by one, Excel
#    rewrites all the dependent formulas automatically.
# ---------------------------------------------------------------------
$wsVar.Rows("11:11").Insert()

# ---------------------------------------------------------------------
# 3) Variables: append the "Signal" piloting block (rows 31-40),
#    mirroring the existing "Aiguillages" block (rows 21-30).
# ---------------------------------------------------------------------
$wsVar.Range("A31").Value = "Signal"

for ($i = 1; $i -le 8; $i++) {
    $r = 31 + $i
    $paramRow = 25 + $i
    $wsVar.Range("A$r").Value = $i
    $wsVar.Range("B$r").Formula = '=$B$11&Variables!A' + $r + '&"="&Paramétrage!B' + $paramRow + '&"&"&"spl"&A' + $r + '&"="&Paramétrage!C' + $paramRow + '&"&spt"&A' + $r + '&"="&Paramétrage!D' + $paramRow + '&"&"'
}

$wsVar.Range("A40").Value = "Total Signal"
$wsVar.Range("B40").Formula = "=B32&B33&B34&B35&B36&B37&B38&B39"

# Now fill in the new header row (done last so "Name Signal"/"sna" get
# appended last in the shared-string table, matching authoring order).
$wsVar.Range("A11").Value = "Name Signal"
$wsVar.Range("B11").Value = "sna"

# ---------------------------------------------------------------------
# 4) Wire the new Signal total into the master config-string formula.
# ---------------------------------------------------------------------
$wsVar.Range("B13").Formula = '=B14&B20&B30&B40&Variables!B9&"="&Paramétrage!B8'

# ---------------------------------------------------------------------
# 5) Paramétrage!B1 now pulls the full string (with signal info) from
#    Variables!B13 (was B12 before the row insert).
# ---------------------------------------------------------------------
$wsParam.Range("B1").Formula = "=Variables!B13"

# ---------------------------------------------------------------------
# 6) Restore the view/selection state recorded in the saved file.
# ---------------------------------------------------------------------
$wsParam.Range("C23").Select()
$wsVar.Range("B30").Select()
$wsVar.Application.ActiveWindow.ScrollRow = 13
